$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339, shifting existing rows 339:380 down to 340:381.
$ws.Rows(339).Insert()

# Populate the newly inserted row 339 with the new data record.
$ws.Range("A339").Value = 10
$ws.Range("B339").Value = "Vega Modelo de Temuco"
$ws.Range("C339").Value = "La Araucanía"
$ws.Range("D339").Value = 45077
$ws.Range("E339").Value = 9
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100103
$ws.Range("H339").Value = "Frutos de hueso (carozo)"
$ws.Range("I339").Value = 100103002
$ws.Range("J339").Value = "Ciruela"
$ws.Range("K339").Value = "Blue Giant"
$ws.Range("L339").Value = "Primera"
$ws.Range("M339").Value = 170
$ws.Range("N339").Value = 14000
$ws.Range("O339").Value = 15000
$ws.Range("P339").Value = 14471
$ws.Range("Q339").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R339").Value = "Región de O'Higgins"
$ws.Range("S339").Value = 804
$ws.Range("T339").Value = 18
